# Update "paises.xlsx" - refresh the timestamp string and a handful of
# per-country statistics (new cases / active cases / recovered / deaths
# columns) to reflect the latest data pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp banner in A1.
$ws.Range("A1").Value = "Datos actualizados a 13 de Octubre de 2020 a las 06:58"

# Row 5 - India (B: total cases, C: new cases, D: recovered, E: active cases)
$ws.Range("B5").Value = 7175880
$ws.Range("C5").Value = 2315
$ws.Range("D5").Value = 6227295
$ws.Range("E5").Value = 838691

# Row 25 - Alemania
$ws.Range("D25").Value = 279100
$ws.Range("E25").Value = 42273

# Row 60 - Uzbekistan
$ws.Range("B60").Value = 61419
$ws.Range("C60").Value = 100
$ws.Range("D60").Value = 58427
$ws.Range("E60").Value = 2483

# Row 68 - Kirguistan
$ws.Range("B68").Value = 49871
$ws.Range("C68").Value = 343
$ws.Range("D68").Value = 44712
$ws.Range("E68").Value = 4067
$ws.Range("G68").Value = 2
$ws.Range("H68").Value = 1092

# Row 143 - Tailandia
$ws.Range("B143").Value = 3643
$ws.Range("C143").Value = 2
$ws.Range("D143").Value = 3457
$ws.Range("E143").Value = 127

# Row 186 - Mongolia
$ws.Range("B186").Value = 320
$ws.Range("C186").Value = 2
$ws.Range("D186").Value = 311
$ws.Range("E186").Value = 9
